$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string (shared string used by A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 21:28"

# 2. Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 100390
$ws.Range("C4").Value = 14955
$ws.Range("E4").Value = 96382
$ws.Range("G4").Value = 248
$ws.Range("H4").Value = 1543

# 3. Alemania (row 8) - updated totals
$ws.Range("B8").Value = 50871
$ws.Range("C8").Value = 6933
$ws.Range("E8").Value = 43871
$ws.Range("G8").Value = 75
$ws.Range("H8").Value = 342

# 4. Jordania (row 78) - updated totals
$ws.Range("D78").Value = 18
$ws.Range("E78").Value = 217

# 5. Insert Montenegro before Martinica (rows 110-112 shuffle down)
# New row 110: Montenegro (fresh data)
$ws.Range("A110").Value = "Montenegro"
$ws.Range("B110").Value = 82
$ws.Range("C110").Value = 13
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 81
$ws.Range("F110").Value = 1
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 1

# New row 111: Martinica (shifted down from old row 110)
$ws.Range("A111").Value = "Martinica"
$ws.Range("B111").Value = 81
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 0
$ws.Range("E111").Value = 80
$ws.Range("F111").Value = 12
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 1

# New row 112: Cuba (shifted down from old row 111)
$ws.Range("A112").Value = "Cuba"
$ws.Range("B112").Value = 80
$ws.Range("C112").Value = 13
$ws.Range("D112").Value = 4
$ws.Range("E112").Value = 74
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 2
